$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# 1) Cell values
# ============================================================

# -- Header row --
$ws.Range("C1").Value = "Session Source"
$ws.Range("D1").Value = "Reccomendation Source"
$ws.Range("G1").Value = "Feedback"

# -- Row 2 (existing row: Flying Lotus) --
$ws.Range("E2").Value = "I really like this artist so far. The sounds were pretty jazzy. A few of the songs had guest artists, which included Kendrick Lamar and Anderson Paak, both of whom I enjoy."
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 1

# -- Row 3: Alice Glass --
$ws.Range("A3").Value = "Alice Glass"
$ws.Range("B3").Value = 44500
$ws.Range("C3").Value = "Amazon Music Artist Channel"
$ws.Range("D3").Value = "Artist similarity to Crystal Castles"
$ws.Range("E3").Value = "Was a former member of Crystal Castles. Generally interesting sounding, but the singers sounded a bit odd to me. I don't remember that voice grating me in Crystal Castles. I think it was similar, but it could be that the songs I listed to didn't break up the singing as much with other parts?"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 1

# -- Row 4: Arkaea --
$ws.Range("A4").Value = "Arkaea"
$ws.Range("B4").Value = 44500
$ws.Range("C4").Value = "You Tube"
$ws.Range("D4").Value = "Close match to artest embedding means for 'The Cure', 'Animals as Leaders', 'System of a Down'"
$ws.Range("E4").Value = "It seems listenable, but makes me want to hear the deftones instead."
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = -1

# -- Row 5: Gorgon City --
$ws.Range("A5").Value = "Gorgon City"
$ws.Range("B5").Value = 44500
$ws.Range("C5").Value = "You Tube"
$ws.Range("D5").Value = "Sample from artist cluster (23): [Electronic, Drum and bass, Trance, EDM, House]"
$ws.Range("E5").Value = "Easy to listen to, but a bit boring."
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = -1

# -- Row 6: Pura Fe --
$ws.Range("A6").Value = "Pura Fé"
$ws.Range("B6").Value = 44500
$ws.Range("C6").Value = "You Tube"
$ws.Range("D6").Value = "Sample from artist cluster (58): [Jazz, Gypsy punk, Pop, Rock, Electric blues, …"
$ws.Range("E6").Value = "The first few songs are pretty interesting. It is different from what I typically listen to. The sounds are folky. It looks like a native American band. "
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 1

# -- Row 7: The Gathering --
$ws.Range("A7").Value = "The Gathering"
$ws.Range("B7").Value = 44500
$ws.Range("C7").Value = "You Tube"
$ws.Range("D7").Value = "Sample from artist cluster (80): [Symphonic metal, Progressive rock, Symphonic …"
$ws.Range("E7").Value = "I really like them so far. It reminds me a bit of night wish, but with a bit softer or ambient sound. There were some sections of songs that came in very different from most of the songs sound, which was nice. The songs feel like they have a musical story/arc. The current song I am on ""Heroes for Ghosts"" has a horn part - awesome!"
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = 1

# -- Row 8: Max Richter --
$ws.Range("A8").Value = "Max Richter"
$ws.Range("B8").Value = 44500
$ws.Range("C8").Value = "Amazon Best of"
$ws.Range("D8").Value = "Sample from artist cluster (17): [Electronic, Ambient, Art rock, Chamber jazz, …"
$ws.Range("E8").Value = "I like the music, but it is something I would want to listen to when going to sleep or concentraing on something else."
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 1

# -- Row 9: Melanie Safka (no Session Source) --
$ws.Range("A9").Value = "Melanie Safka"
$ws.Range("B9").Value = 44500
$ws.Range("D9").Value = "Sample from artist cluster (7): [Alternative rock, Trip hop, post-industrial, …"
$ws.Range("E9").Value = "I didn't expect her sound based on the top tagged genres. I like the songs so far. I recognized some of the songs by her and some sounded like covers. The sound seems more 70's folk style."
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = 1

# ============================================================
# 2) Formatting - copy styles from existing reference cells so
#    that shared style indices are reused instead of minted fresh
# ============================================================

# style "2" (wrap + vertical-top) reference = A2
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:E3").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:E5").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C6:E6").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7:E8").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D9:E9").PasteSpecial(-4122) | Out-Null
# note: A6 & A9 intentionally keep the default (no) style, matching source

# style "1" (vertical-top only) reference = F2
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F3:F9").PasteSpecial(-4122) | Out-Null

# style "3" (numFmtId 14, date) reference = B2
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4:B9").PasteSpecial(-4122) | Out-Null

# style "4" (numFmtId 16, date w/o year) - new format, only B3 uses it
$ws.Range("B3").NumberFormat = "d-mmm"

# A6 & A9 revert to the plain default/"Normal" style (no wrap/valign),
# overriding the column-level style they would otherwise inherit
$ws.Range("A6").Style = "Normal"
$ws.Range("A9").Style = "Normal"

# ============================================================
# 3) Column widths / row heights
# ============================================================
$ws.Columns.Item(1).ColumnWidth = 12.54296875
$ws.Columns.Item(4).ColumnWidth = 32.26953125
$ws.Columns.Item(5).ColumnWidth = 38.7265625

$ws.Rows.Item(3).RowHeight = 101.5
$ws.Rows.Item(4).RowHeight = 43.5
$ws.Rows.Item(5).RowHeight = 43.5
$ws.Rows.Item(6).RowHeight = 58
$ws.Rows.Item(7).RowHeight = 116
$ws.Rows.Item(8).RowHeight = 43.5
$ws.Rows.Item(9).RowHeight = 72.5

# ============================================================
# 4) View state
# ============================================================
$ws.Range("F7").Select() | Out-Null

Write-Output "done"
